# Actualizacion 11 de Mayo - Tarde
# Updates the "6ARHM" group statistics on the 2P / Final summary sheets and
# refreshes the "Rescatables" (pending/failing students) roster: two new
# students are added and one student's failing-subject count drops from 2 to 1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Estadisticos 2P" - row 4 is the 6ARHM group.
# ---------------------------------------------------------------------------
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")
$ws2P.Range("D4").Value = 8
$ws2P.Range("E4").Value = 10
$ws2P.Range("F4").Value = 24
$ws2P.Range("G4").Value = 70.59
$ws2P.Range("H4").Value = 8.300000000000001

# ---------------------------------------------------------------------------
# 2) "Estadisticos Final" - row 4 is the 6ARHM group.
# ---------------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Range("E4").Value = 6
$wsFinal.Range("F4").Value = 28
$wsFinal.Range("G4").Value = 82.34999999999999

# ---------------------------------------------------------------------------
# 3) "Rescatables" - full roster rewrite (rows 2:23), reflecting:
#      - a new student (SANTIAGO CRUZ YARITZI) inserted at the top
#      - a new student (ROJAS MAZA ANGEL GABRIEL) inserted
#      - RUIZ LOPEZ ALFONSO now only fails 1 subject (was 2) and moves down
#      - two new students appended at the bottom (MAZAHUA IXMATLAHUA SOFIA,
#        ZARATE CASTILLO MIGUEL ANGEL)
# ---------------------------------------------------------------------------
$wsResc = $wb.Worksheets.Item("Rescatables")

$probEst = "PROBABILIDAD Y ESTADÍSTICA"

$data = @(
    @(19330051920117, "SANTIAGO",  "CRUZ",       "YARITZI",       "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "4BEM",  2),
    @(18330051920152, "CRISTOBAL", "ROMERO",     "EDGAR ARMANDO", $probEst, "6AEM",  2),
    @(18330051920160, "MARRON",    "BLASQUEZ",   "DAVID",         $probEst, "6AEM",  2),
    @(18330051920172, "RIOS",      "OCHOA",      "JONATHAN",      $probEst, "6AEM",  2),
    @(18330051920176, "SANCHEZ",   "TRUJILLO",   "ERIK JAIR",     $probEst, "6AEM",  2),
    @(18330051920026, "ROJAS",     "MAZA",       "ANGEL GABRIEL", $probEst, "6AEV",  2),
    @(18330051920217, "MARQUEZ",   "HERNANDEZ",  "ANDRES",        $probEst, "6ARHM", 2),
    @(18330051920329, "GUERRA",    "ROMERO",     "JOCELYN",       $probEst, "6BLCM", 2),
    @(18330051920424, "LOPEZ",     "APALE",      "MARIA LARET",   $probEst, "6BLCM", 2),
    @(18330051920339, "MORALES",   "TREJO",      "ROCIO TAMARA",  $probEst, "6BLCM", 2),
    @(18330051920342, "REYES",     "SARMIENTO",  "INGRID PAOLA",  $probEst, "6BLCM", 2),
    @(18330051920347, "TEXOCO",    "DE JESUS",   "MAYTE",         $probEst, "6BLCM", 2),
    @(18330051920352, "ZEPAHUA",   "JUAREZ",     "QUETZALI",      $probEst, "6BLCM", 2),
    @(18330051920143, "ARIAS",     "HERNANDEZ",  "JUAN ALBERTO",  $probEst, "6AEM",  1),
    @(18330051920141, "ANTONIO",   "OFICIAL",    "EDUARDO",       $probEst, "6AEM",  1),
    @(17330051920160, "RUIZ",      "LOPEZ",      "ALFONSO",       $probEst, "6AEV",  1),
    @(18330051920094, "DIAZ",      "SESMA",      "FLOR LIZZETH",  $probEst, "6ARHM", 1),
    @(18330051920223, "MIXCOHUA",  "ZEPAHUA",    "REBECA MELISA", $probEst, "6ARHM", 1),
    @(18330051920412, "SANCHEZ",   "SORIANO",    "ANTONIO",       $probEst, "6ASM",  1),
    @(18330051920334, "LOPEZ",     "ROJAS",      "JOEL",          $probEst, "6BLCM", 1),
    @(18330051920337, "MAZAHUA",   "IXMATLAHUA", "SOFIA",         $probEst, "6BLCM", 1),
    @(18330051920427, "ZARATE",    "CASTILLO",   "MIGUEL ANGEL",  $probEst, "6BLCM", 1)
)

$r = 2
foreach ($row in $data) {
    $wsResc.Cells.Item($r, 1).Value = $row[0]
    $wsResc.Cells.Item($r, 2).Value = $row[1]
    $wsResc.Cells.Item($r, 3).Value = $row[2]
    $wsResc.Cells.Item($r, 4).Value = $row[3]
    $wsResc.Cells.Item($r, 5).Value = $row[4]
    $wsResc.Cells.Item($r, 6).Value = $row[5]
    $wsResc.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
